$d = $word.ActiveDocument
$d.Save()
